$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "I can open the file from Github--but  I can't see the changes I made in the file"
